# Auto-generated: applies scheduled-runner market/profit data refresh
# to the Leviathan_Profits workbook (columns H-N per leve row) across all
# eight job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 2075
$ws.Range("J19").Value = 610
$ws.Range("L19").Value = 610
$ws.Range("N19").Value = -960
# row 43
$ws.Range("H43").Value = 3332.7693
$ws.Range("I43").Value = 2228.3333
$ws.Range("J43").Value = 4279.4287
$ws.Range("K43").Value = 2228.3333
$ws.Range("L43").Value = 4279.4287
$ws.Range("M43").Value = -2159.3333
$ws.Range("N43").Value = -4417.4287
# row 61
$ws.Range("H61").Value = 1002.6667
$ws.Range("I61").Value = 999.8
$ws.Range("J61").Value = 1017
$ws.Range("K61").Value = 2999.4
$ws.Range("L61").Value = 3051
$ws.Range("M61").Value = -2827.4
$ws.Range("N61").Value = -3395
# row 86
$ws.Range("H86").Value = 1741.2667
$ws.Range("I86").Value = 1429.3334
$ws.Range("J86").Value = 1949.2222
$ws.Range("K86").Value = 1429.3334
$ws.Range("L86").Value = 1949.2222
$ws.Range("M86").Value = -306.3334
$ws.Range("N86").Value = -4195.2222
# row 89
$ws.Range("H89").Value = 1741.2667
$ws.Range("I89").Value = 1429.3334
$ws.Range("J89").Value = 1949.2222
$ws.Range("K89").Value = 7146.666999999999
$ws.Range("L89").Value = 9746.110999999999
$ws.Range("M89").Value = -1530.666999999999
$ws.Range("N89").Value = -20978.111
# row 98
$ws.Range("H98").Value = 12481.277
$ws.Range("I98").Value = 12566.9
$ws.Range("J98").Value = 12374.25
$ws.Range("K98").Value = 12566.9
$ws.Range("L98").Value = 12374.25
$ws.Range("M98").Value = -11068.9
$ws.Range("N98").Value = -15370.25
# row 116
$ws.Range("H116").Value = 39999.5
$ws.Range("J116").Value = 39999.5
$ws.Range("L116").Value = 39999.5
$ws.Range("N116").Value = -46883.5
# row 122
$ws.Range("H122").Value = 12481.277
$ws.Range("I122").Value = 12566.9
$ws.Range("J122").Value = 12374.25
$ws.Range("K122").Value = 37700.7
$ws.Range("L122").Value = 37122.75
$ws.Range("M122").Value = -35250.7
$ws.Range("N122").Value = -42022.75
# row 131
$ws.Range("H131").Value = 646.5
$ws.Range("I131").Value = 1000
$ws.Range("K131").Value = 3000
$ws.Range("M131").Value = 2040
# row 132
$ws.Range("H132").Value = 3566.6667
$ws.Range("I132").Value = 1347.7667
$ws.Range("J132").Value = 10963
$ws.Range("K132").Value = 4043.300099999999
$ws.Range("L132").Value = 32889
$ws.Range("M132").Value = -1513.300099999999
$ws.Range("N132").Value = -37949
# row 137
$ws.Range("H137").Value = 45637.086
$ws.Range("I137").Value = 2145
$ws.Range("J137").Value = 252224.5
$ws.Range("K137").Value = 6435
$ws.Range("L137").Value = 756673.5
$ws.Range("M137").Value = -3885
$ws.Range("N137").Value = -761773.5
# row 138
$ws.Range("H138").Value = 1808.1364
$ws.Range("I138").Value = 1230.4584
$ws.Range("K138").Value = 3691.3752
$ws.Range("M138").Value = 1448.6248

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 11199.88
$ws.Range("I45").Value = 8653.294
$ws.Range("J45").Value = 16611.375
$ws.Range("K45").Value = 8653.294
$ws.Range("L45").Value = 16611.375
$ws.Range("M45").Value = -8276.294
$ws.Range("N45").Value = -17365.375
# row 61
$ws.Range("H61").Value = 2777.1765
$ws.Range("I61").Value = 2747.5334
$ws.Range("K61").Value = 2747.5334
$ws.Range("M61").Value = -2535.5334
# row 136
$ws.Range("H136").Value = 2777.1765
$ws.Range("I136").Value = 2747.5334
$ws.Range("K136").Value = 8242.600199999999
$ws.Range("M136").Value = -5692.600199999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 80
$ws.Range("H80").Value = 798.8889
$ws.Range("I80").Value = 660.5
$ws.Range("J80").Value = 880.2941
$ws.Range("K80").Value = 660.5
$ws.Range("L80").Value = 880.2941
$ws.Range("M80").Value = 337.5
$ws.Range("N80").Value = -2876.2941
# row 83
$ws.Range("H83").Value = 798.8889
$ws.Range("I83").Value = 660.5
$ws.Range("J83").Value = 880.2941
$ws.Range("K83").Value = 3302.5
$ws.Range("L83").Value = 4401.470499999999
$ws.Range("M83").Value = 1689.5
$ws.Range("N83").Value = -14385.4705
# row 86
$ws.Range("H86").Value = 1291.119
$ws.Range("I86").Value = 1076.9412
$ws.Range("J86").Value = 2201.375
$ws.Range("K86").Value = 1076.9412
$ws.Range("L86").Value = 2201.375
$ws.Range("M86").Value = 46.05880000000002
$ws.Range("N86").Value = -4447.375
# row 89
$ws.Range("H89").Value = 1291.119
$ws.Range("I89").Value = 1076.9412
$ws.Range("J89").Value = 2201.375
$ws.Range("K89").Value = 5384.706
$ws.Range("L89").Value = 11006.875
$ws.Range("M89").Value = 231.2939999999999
$ws.Range("N89").Value = -22238.875
# row 134
$ws.Range("H134").Value = 2446.9473
$ws.Range("I134").Value = 2128.0625
$ws.Range("K134").Value = 6384.1875
$ws.Range("M134").Value = -3849.1875

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2882.1904
$ws.Range("I31").Value = 1466.1875
$ws.Range("J31").Value = 7413.4
$ws.Range("K31").Value = 1466.1875
$ws.Range("L31").Value = 7413.4
$ws.Range("M31").Value = -1171.1875
$ws.Range("N31").Value = -8003.4
# row 34
$ws.Range("H34").Value = 2882.1904
$ws.Range("I34").Value = 1466.1875
$ws.Range("J34").Value = 7413.4
$ws.Range("K34").Value = 1466.1875
$ws.Range("L34").Value = 7413.4
$ws.Range("M34").Value = -1264.1875
$ws.Range("N34").Value = -7817.4
# row 58
$ws.Range("H58").Value = 1608.7858
$ws.Range("J58").Value = 2046.2
$ws.Range("L58").Value = 2046.2
$ws.Range("N58").Value = -2452.2
# row 105
$ws.Range("H105").Value = 3261.75
$ws.Range("I105").Value = 3544.7222
$ws.Range("J105").Value = 2978.7778
$ws.Range("K105").Value = 3544.7222
$ws.Range("L105").Value = 2978.7778
$ws.Range("M105").Value = -1797.7222
$ws.Range("N105").Value = -6472.7778
# row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
# row 134
$ws.Range("H134").Value = 1937.4231
$ws.Range("I134").Value = 1710.9524
$ws.Range("K134").Value = 5132.857199999999
$ws.Range("M134").Value = -2597.857199999999
# row 136
$ws.Range("H136").Value = 1608.7858
$ws.Range("J136").Value = 2046.2
$ws.Range("L136").Value = 6138.6
$ws.Range("N136").Value = -11238.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 23
$ws.Range("H23").Value = 831.8570999999999
$ws.Range("J23").Value = 867.4
$ws.Range("L23").Value = 2602.2
$ws.Range("N23").Value = -3072.2
# row 52
$ws.Range("H52").Value = 823
$ws.Range("J52").Value = 823
$ws.Range("L52").Value = 2469
$ws.Range("N52").Value = -3001
# row 80
$ws.Range("J80").Value = 2333.3333
$ws.Range("L80").Value = 6999.999899999999
$ws.Range("N80").Value = -8871.999899999999
# row 83
$ws.Range("J83").Value = 2333.3333
$ws.Range("L83").Value = 20999.9997
$ws.Range("N83").Value = -30359.9997
# row 107
$ws.Range("H107").Value = 1341.4706
$ws.Range("I107").Value = 1500
$ws.Range("J107").Value = 1320.3334
$ws.Range("K107").Value = 4500
$ws.Range("L107").Value = 3961.0002
$ws.Range("M107").Value = -2580
$ws.Range("N107").Value = -7801.0002
# row 140
$ws.Range("H140").Value = 3606.8
$ws.Range("I140").Value = 2587.2
$ws.Range("K140").Value = 7761.599999999999
$ws.Range("M140").Value = -2581.599999999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 11768427
$ws.Range("I70").Value = 13336897
$ws.Range("J70").Value = 4899.5
$ws.Range("K70").Value = 13336897
$ws.Range("L70").Value = 4899.5
$ws.Range("M70").Value = -13336627
$ws.Range("N70").Value = -5439.5
# row 73
$ws.Range("H73").Value = 11768427
$ws.Range("I73").Value = 13336897
$ws.Range("J73").Value = 4899.5
$ws.Range("K73").Value = 13336897
$ws.Range("L73").Value = 4899.5
$ws.Range("M73").Value = -13335961
$ws.Range("N73").Value = -6771.5
# row 122
$ws.Range("H122").Value = 2499.5
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -13897
# row 123
$ws.Range("H123").Value = 65217
$ws.Range("J123").Value = 65217
$ws.Range("L123").Value = 65217
$ws.Range("N123").Value = -70117

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 5947.4
$ws.Range("I40").Value = 5182.857
$ws.Range("K40").Value = 5182.857
$ws.Range("M40").Value = -5046.857

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# row 132
$ws.Range("H132").Value = 9452.261
$ws.Range("I132").Value = 10662.723
$ws.Range("K132").Value = 31988.169
$ws.Range("M132").Value = -29458.169
